# repull data, push all data, mean calculation
# Update the dSF (column F) values for the affected rows to reflect
# the re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 3
    12 = 0
    30 = 6
    31 = -5
    38 = 1
    39 = 1
    40 = -2
    44 = 0
    47 = 3
    52 = 9
    57 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
